# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the "Estado de Cuenta" worksheet with the new worker's data and
# mora (overdue) amounts, and widens column D to fit the longer name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Valor Mora (header summary box) - matches the new detail row below
$ws.Range("E11").Value = 40000

# Worker detail row (row 16): document type, document number, name, period
$ws.Range("B16").Value = "PE"
$ws.Range("C16").Value = "810938027101980"
$ws.Range("D16").Value = "EDGAR ALEXANDER RAMIREZ BELANDRIA"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1000000

# Widen the "Nombre Trabajador" column so the longer name still fits
# (ColumnWidth is in "characters"; Excel stores the XML width in a slightly
# different unit, so back the character width off by the fixed 5/6 padding
# offset this workbook's font applies, landing exactly on width=38).
$ws.Columns("D").ColumnWidth = (38 - 5/6)
